# The "manager" added a new menu item: row 5 (Item ID 4) is renamed
# from "test" to "Pizza". The price (1), Special Item (false) and
# Is Active (false) flags for that row stay the same.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "Pizza"
